$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 999.8333
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H56").Value = 999.8333
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H58").Value = 125
$ws.Range("J58").Value = 110
$ws.Range("L58").Value = 330
$ws.Range("N58").Value = -630
$ws.Range("H62").Value = 6182
$ws.Range("I62").Value = 5546.077
$ws.Range("K62").Value = 5546.077
$ws.Range("M62").Value = -4922.077
$ws.Range("H65").Value = 6182
$ws.Range("I65").Value = 5546.077
$ws.Range("K65").Value = 27730.385
$ws.Range("M65").Value = -24610.385
$ws.Range("H76").Value = 5000
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 5000
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H135").Value = 506
$ws.Range("I135").Value = 506
$ws.Range("K135").Value = 4554
$ws.Range("M135").Value = -2019

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 121.6
$ws.Range("I4").Value = 121.6
$ws.Range("K4").Value = 121.6
$ws.Range("M4").Value = -5.599999999999994
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
$ws.Range("H110").Value = 796
$ws.Range("I110").Value = 728
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 728
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1317
$ws.Range("N110").Value = -5090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 2263.8333
$ws.Range("I25").Value = 2016.6
$ws.Range("K25").Value = 2016.6
$ws.Range("M25").Value = -1781.6
$ws.Range("H64").Value = 688.5
$ws.Range("I64").Value = 688.5
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 688.5
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -463.5
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 688.5
$ws.Range("I67").Value = 688.5
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 688.5
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = 91.5
$ws.Range("N67").ClearContents()
$ws.Range("H94").Value = 631.8570999999999
$ws.Range("I94").Value = 631.8570999999999
$ws.Range("K94").Value = 631.8570999999999
$ws.Range("M94").Value = -180.8570999999999
$ws.Range("H107").Value = 1611
$ws.Range("I107").Value = 1611
$ws.Range("K107").Value = 1611
$ws.Range("M107").Value = 309

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2699.6667
$ws.Range("I62").Value = 2699.6667
$ws.Range("K62").Value = 2699.6667
$ws.Range("M62").Value = -2075.6667
$ws.Range("H65").Value = 2699.6667
$ws.Range("I65").Value = 2699.6667
$ws.Range("K65").Value = 13498.3335
$ws.Range("M65").Value = -10378.3335
$ws.Range("H94").Value = 849.5
$ws.Range("I94").Value = 849
$ws.Range("K94").Value = 849
$ws.Range("M94").Value = -398
$ws.Range("H107").Value = 1138.75
$ws.Range("I107").Value = 1555.5
$ws.Range("J107").Value = 999.8333
$ws.Range("K107").Value = 1555.5
$ws.Range("L107").Value = 999.8333
$ws.Range("M107").Value = 364.5
$ws.Range("N107").Value = -4839.8333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4393.5
$ws.Range("I63").Value = 4393.5
$ws.Range("K63").Value = 13180.5
$ws.Range("M63").Value = -12431.5
$ws.Range("H66").Value = 4393.5
$ws.Range("I66").Value = 4393.5
$ws.Range("K66").Value = 39541.5
$ws.Range("M66").Value = -35797.5
$ws.Range("H98").Value = 1445
$ws.Range("I98").Value = 1445
$ws.Range("K98").Value = 4335
$ws.Range("M98").Value = -2837

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H99").Value = 14666.667
$ws.Range("I99").Value = 14666.667
$ws.Range("K99").Value = 14666.667
$ws.Range("M99").Value = -12420.667
$ws.Range("H107").Value = 410
$ws.Range("I107").Value = 350
$ws.Range("J107").Value = 470
$ws.Range("K107").Value = 350
$ws.Range("L107").Value = 470
$ws.Range("M107").Value = 1570
$ws.Range("N107").Value = -4310
$ws.Range("H122").Value = 999.5
$ws.Range("I122").Value = 999.5
$ws.Range("K122").Value = 2998.5
$ws.Range("M122").Value = -548.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 199
$ws.Range("I40").Value = 199
$ws.Range("K40").Value = 199
$ws.Range("M40").Value = -63
$ws.Range("H55").Value = 440
$ws.Range("I55").Value = 420
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 420
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = -247
$ws.Range("N55").Value = -846
$ws.Range("H122").Value = 1662.6666
$ws.Range("I122").Value = 1662.6666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4987.9998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2537.9998
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5216.5
$ws.Range("I62").Value = 4759.8
$ws.Range("K62").Value = 4759.8
$ws.Range("M62").Value = -4135.8
$ws.Range("H65").Value = 5216.5
$ws.Range("I65").Value = 4759.8
$ws.Range("K65").Value = 23799
$ws.Range("M65").Value = -20679
$ws.Range("H96").Value = 1850
$ws.Range("I96").Value = 1966.6666
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 1966.6666
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = -593.6666
$ws.Range("N96").Value = -4246
